# Bug fix in Eduati data files:
#  - Sheet1 ("CCK81_noCTRL_meas") had 43 stray extra rows (A45:A87, a leftover
#    sequential-number fill) below the real 44-row data table; trim them off.
#  - Restore the view/selection state that Excel wrote when the file was last
#    saved: Sheet1 becomes the active/selected tab (cell C68 selected,
#    scrolled to row 24), and Sheet3 is no longer the selected tab.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")

# Sheet1 dimension was A1:N87 but only A1:N44 holds real data -- rows 45:87
# are a stray leftover (column A only, values 44..86) that should not be
# there. Delete them so the sheet shrinks back down to A1:N44.
$ws1.Rows("45:87").Delete()

# Make Sheet1 the active sheet/tab (was Sheet3), move the view to where it
# was left (scrolled down so row 24 is at the top) and select C68.
$ws1.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 24
$win.ScrollColumn = 1
$ws1.Range("C68").Select()
